$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: text labels (rows 2-8) - write first so shared strings are allocated in this order
$ws.Range("A2").Value = "someone "
$ws.Range("A3").Value = "ape-thief"
$ws.Range("A4").Value = "woman"
$ws.Range("A5").Value = " gorilla"
$ws.Range("A6").Value = "that"
$ws.Range("A7").Value = "batman"
$ws.Range("A8").Value = "someone"

# Column B: absolute-path centric labels (rows 2-8) - write second
$ws.Range("B2").Value = "root_Entity_Per_Per"
$ws.Range("B3").Value = "root_Event_Life_Be-Born_Trigger"
$ws.Range("B4").Value = "root_Event_Life_Marry_Trigger"
$ws.Range("B5").Value = "root_Event_Life_Marry_Person"
$ws.Range("B6").Value = "root_Event_Life_Marry_Place"
$ws.Range("B7").Value = "root_Entity_Veh_Veh"
$ws.Range("B8").Value = "root_Entity_Veh_Veh"

# Column C: counts (rows 2-8) all become 1
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 1
$ws.Range("C4").Value = 1
$ws.Range("C5").Value = 1
$ws.Range("C6").Value = 1
$ws.Range("C7").Value = 1
$ws.Range("C8").Value = 1
